$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 47, shifting rows 47:218 down to 48:219
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47
$ws.Cells.Item(47, 1).Value = "area"
$ws.Cells.Item(47, 2).Value = "QUALITY"
$ws.Cells.Item(47, 3).Value = 1
$ws.Cells.Item(47, 4).Value = 7
